$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (tipo-nacionalidad)
$ws.Range("A2").Value = "iaest-measure:tipo-nacionalidad"
$ws.Range("A3").Value = "medida"
$ws.Range("A4").Value = "xsd:int"

# Column C (continente-nacionalidad)
$ws.Range("C2").Value = "iaest-measure:continente-nacionalidad"
$ws.Range("C3").Value = "medida"
$ws.Range("C4").Value = "xsd:int"

# Column E (municipio-nombre)
$ws.Range("E2").Value = "sdmx-dimension:refArea"
$ws.Range("E3").Value = "dim"
$ws.Range("E4").Value = "URI-Municipio"

# Column H (aragon)
$ws.Range("H2").Value = "sdmx-dimension:refArea"
$ws.Range("H4").Value = "URI-Comunidad"

# Row 5 (mapping file references) is removed entirely
$ws.Rows("5:5").Delete()
